# "Expanded with comparison of random strings" -------------------------------
# Adds a second table comparing random strings (decimal/hex digits, letters,
# mixed-case letters, binary) against the word-length table above it, plus a
# couple of annotation cells, and widens column A to fit the longer labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ------------------------------------------------------------------------
# 0. Preserve the two "source for ..." footnote rows before we overwrite
#    the cells they currently live in (old rows 13 & 14 -> new rows 20 & 21).
# ------------------------------------------------------------------------
$sourceWordCounts = $ws.Range("B13").Value2
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A20").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B20").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B21").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("B20").Value2 = $sourceWordCounts
$ws.Range("B21").Value2 = "source for Dicewords:  http://world.std.com/~reinhold/diceware.html"

# Clear the old B14 (Dicewords footnote) text now that it has moved to B21,
# and the old B13 text now that it has moved to B20; both cells get new
# content below.
$ws.Range("B13").ClearContents() | Out-Null
$ws.Range("B14").ClearContents() | Out-Null

# ------------------------------------------------------------------------
# 1. Column A width so the new, longer row labels fit (target ~31.57
#    characters; 30.67 is the input that round-trips closest to that
#    stored width given this engine's column-width quantization).
# ------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 30.67

# ------------------------------------------------------------------------
# 2. Annotation cells next to the existing word-length table.
# ------------------------------------------------------------------------
$ws.Range("A1").Copy() | Out-Null
$ws.Range("J9").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("J11").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("J9").Value2 = "The clear winner"
$ws.Range("J11").Value2 = "(6^5 words, 32961 characters total)"

# ------------------------------------------------------------------------
# 3. Row 12: extend the blank spacer row's formatting to C12/D12.
# ------------------------------------------------------------------------
$ws.Range("A1").Copy() | Out-Null
$ws.Range("C12").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D12").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

# ------------------------------------------------------------------------
# 4. Row 13: new section header.
# ------------------------------------------------------------------------
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A13").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B13").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C13").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D13").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A13").Value2 = "Random strings of characters, for comparison"

# ------------------------------------------------------------------------
# 5. Rows 14-18: the new random-strings comparison table, mirroring the
#    word-length table's layout/formulas exactly.
# ------------------------------------------------------------------------
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A14:C18").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D14:D18").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A14").Value2 = "Decimal digits"
$ws.Range("A15").Value2 = "Hex digits"
$ws.Range("A16").Value2 = "Random letters"
$ws.Range("A17").Value2 = "Mixed case letters"
$ws.Range("A18").Value2 = "Binary"

$ws.Range("B14").Value2 = 1.0
$ws.Range("C14").Value2 = 10.0
$ws.Range("B15").Value2 = 1.0
$ws.Range("C15").Value2 = 16.0
$ws.Range("B16").Value2 = 1.0
$ws.Range("C16").Value2 = 26.0
$ws.Range("B17").Value2 = 1.0
$ws.Range("C17").Value2 = 52.0
$ws.Range("B18").Value2 = 1.0
$ws.Range("C18").Value2 = 256.0

$ws.Range("D14:D18").Formula = "=1.44*LN(C14)"
$ws.Range("E14:E18").Formula = "=1.44*LN(C14^2)"
$ws.Range("F14:F18").Formula = "=1.44*LN(C14^3)"
$ws.Range("G14:G18").Formula = "=1.44*LN(C14^4)"
$ws.Range("H14:H18").Formula = "=1.44*LN(C14^5)"
$ws.Range("I14:I18").Formula = "=D14/(B14)"

# Annotation cells for the new table.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("J15").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("J18").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("J15").Value2 = "(should be exactly 4.0, rounding errors due to use of natural log)"
$ws.Range("J18").Value2 = "(same as above, but 8.0)"

# ------------------------------------------------------------------------
# 6. Row 19: blank spacer row below the new table.
# ------------------------------------------------------------------------
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A19").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B19").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false
